# CIERRE 19 MAR 22
# Advance the payroll week: "SEMANA 10 ... 07 al 13 MARZO 2022" -> "SEMANA 11 ... 14 al 20 MARZO 2022"
# Rotate the employee "# n" turn markers and record this week's EXTRAS/PRESTAMO payment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Week header text (B9 is the source cell; H9/B27/H27/B43/H43 are formulas that
#        reference it, so they recalc automatically once B9's value changes). ---
$ws.Range("B9").Value = "SEMANA   11  DEL    14      Al   20   DE   MARZO          2022"

# --- 2. Rotate the "# n" employee-turn labels. ---
$ws.Range("H14").Value = "# 4"
$ws.Range("H32").Value = "# 5"
$ws.Range("B48").Value = "# 3"
# H48 ("# 6") is unchanged this week.

# --- 3. This period's EXTRAS / PRESTAMO amount for TEODORA ARELLANO PEREZ (row 38-41). ---
$ws.Range("K39").ClearContents()
$ws.Range("K40").Value = 1250

# --- 4. Update the view's selection to match where the author left off. ---
$ws.Range("H33").Select() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 10
$aw.ScrollColumn = 1
